# [Justice Counts][2/n] Get rid of "All" value in Tech Spec
#
# 1. Rename the "total_staff" sheet to "total_staff_by_type".
# 2. Move the active-cell selection on that sheet from D2 to D17.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("total_staff")
$ws.Name = "total_staff_by_type"

# Activate the sheet and move the selection to D17 (mirrors the saved
# <selection activeCell="D17" sqref="D17"/> in the sheet's sheetView).
$ws.Activate()
$ws.Range("D17").Select()
